# Add a new real-time GDP vintage column (CF) that mirrors column CE for
# existing rows, with a couple of revised/new data points at the tail, and
# append one new observation row (141) for the new reference period.
#
# Column CE = 83, Column CF = 84 (1-based column indices)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Header cell CF1: new vintage/release date (copy CE1's date style) ---
$ceCell = $ws.Cells.Item(1, 83)
$cfCell = $ws.Cells.Item(1, 84)
$ceCell.Copy()
$cfCell.PasteSpecial($xlPasteFormats)
$cfCell.Value = 45986

# --- Data rows 2 .. 138: CF = CE (same vintage value carried into new column) ---
for ($r = 2; $r -le 138; $r++) {
    $ceCell = $ws.Cells.Item($r, 83)
    $ceVal = $ceCell.Value2
    if ($ceVal -ne "") {
        $cfCell = $ws.Cells.Item($r, 84)
        $ceCell.Copy()
        $cfCell.PasteSpecial($xlPasteFormats)
        $cfCell.Value = $ceVal
    }
}

# --- Row 139: revised value in new vintage (104.52 -> 104.59) ---
$ceCell139 = $ws.Cells.Item(139, 83)
$cfCell139 = $ws.Cells.Item(139, 84)
$ceCell139.Copy()
$cfCell139.PasteSpecial($xlPasteFormats)
$cfCell139.Value = 104.59

# --- Row 140: first appearance of this period's value in the new vintage ---
$ws.Cells.Item(140, 84).Value = 104.59

# --- Row 141: brand-new observation period, only column A populated ---
$ws.Range("A140").Copy()
$ws.Range("A141").PasteSpecial($xlPasteFormats)
$ws.Range("A141").Value = 45976

$excel.CutCopyMode = $false
